$d = $word.ActiveDocument

# Locate the "KEY ACHIEVEMENTS AND IMPACT" heading paragraph, then operate on
# the bullet paragraphs that follow its "Impact" sub-heading. Doing this by
# text search (rather than a hard-coded index) keeps the script resilient to
# any earlier paragraphs in the document.
$count = $d.Paragraphs.Count
$headingIdx = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*KEY ACHIEVEMENTS AND IMPACT*") {
        $headingIdx = $i
        break
    }
}

if ($headingIdx -eq -1) {
    throw "Could not find 'KEY ACHIEVEMENTS AND IMPACT' heading"
}

# Paragraph layout relative to the heading:
#   headingIdx     -> "KEY ACHIEVEMENTS AND IMPACT"
#   headingIdx + 1 -> "Impact" (Heading3)
#   headingIdx + 2 -> bullet 1 (Achieved 87% prediction accuracy...)
#   headingIdx + 3 -> bullet 2 (Delivered $4.9M additional revenue...)
#   headingIdx + 4 -> bullet 3 (Built redistricting platform...)
#   headingIdx + 5 -> bullet 4 (Developed longitudinal data analysis...)
#   headingIdx + 6 -> bullet 5 (Discovered systematic race coding errors...)
#   headingIdx + 7 -> bullet 6 (Trigonometric algorithm for boundary...)

$bullet1 = $headingIdx + 2
$bullet2 = $headingIdx + 3
$bullet3 = $headingIdx + 4
$bullet4 = $headingIdx + 5
$bullet5 = $headingIdx + 6
$bullet6 = $headingIdx + 7

# Sanity-check the paragraphs we are about to touch contain the text we
# expect, so we fail loudly instead of silently corrupting the wrong
# paragraphs if the document layout ever shifts.
if ($d.Paragraphs.Item($bullet1).Range.Text -notlike "*Achieved 87% prediction accuracy*") {
    throw "Unexpected text at bullet1"
}
if ($d.Paragraphs.Item($bullet2).Range.Text -notlike "*Delivered `$4.9M additional revenue*") {
    throw "Unexpected text at bullet2"
}
if ($d.Paragraphs.Item($bullet3).Range.Text -notlike "*Built redistricting platform*") {
    throw "Unexpected text at bullet3"
}
if ($d.Paragraphs.Item($bullet4).Range.Text -notlike "*Developed longitudinal data analysis*") {
    throw "Unexpected text at bullet4"
}
if ($d.Paragraphs.Item($bullet5).Range.Text -notlike "*Discovered systematic race coding errors*") {
    throw "Unexpected text at bullet5"
}
if ($d.Paragraphs.Item($bullet6).Range.Text -notlike "*Trigonometric algorithm for boundary*") {
    throw "Unexpected text at bullet6"
}

# Rewrite the first four bullets with the new accomplishment-focused text.
$d.Paragraphs.Item($bullet1).Range.Text = "• Revenue generation: Delivered `$4.9M additional revenue through optimization"
$d.Paragraphs.Item($bullet2).Range.Text = "• 23% conversion rate improvement"
$d.Paragraphs.Item($bullet3).Range.Text = "• Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis"
$d.Paragraphs.Item($bullet4).Range.Text = "• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"

# Remove the two trailing bullets entirely (delete from the end so the
# earlier index stays valid while deleting).
$d.Paragraphs.Item($bullet6).Range.Delete()
$d.Paragraphs.Item($bullet5).Range.Delete()
